$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$c = $ws.Range("D2"); $s = $c.Style; $c.NumberFormat = "@"; $c.Value = '60.813.61'; $c.Style = $s
$ws.Range("E2").Value = '  -2.55%  '

# Row 3 - Ethereum
$c = $ws.Range("D3"); $s = $c.Style; $c.NumberFormat = "@"; $c.Value = '2.904.18'; $c.Style = $s
$ws.Range("E3").Value = '  -3.81%  '

# Row 4 - TetherUSD
$c = $ws.Range("D4"); $s = $c.Style; $c.NumberFormat = "@"; $c.Value = '1.00'; $c.Style = $s
$ws.Range("E4").Value = '  -0.09%  '

# Row 5 - BNB
$c = $ws.Range("D5"); $s = $c.Style; $c.NumberFormat = "@"; $c.Value = '586.68'; $c.Style = $s
$ws.Range("E5").Value = '  -1.46%  '

# Row 6 - Solana
$c = $ws.Range("D6"); $s = $c.Style; $c.NumberFormat = "@"; $c.Value = '146.99'; $c.Style = $s
$ws.Range("E6").Value = '  -1.85%  '

# Row 7 - USDC
$ws.Range("E7").Value = '  +0.09%  '

# Row 8 - XRP
$c = $ws.Range("D8"); $s = $c.Style; $c.NumberFormat = "@"; $c.Value = '0.505'; $c.Style = $s
$ws.Range("E8").Value = '  -2.75%  '

# Row 9 - LidoStakedEther
$c = $ws.Range("D9"); $s = $c.Style; $c.NumberFormat = "@"; $c.Value = '2.905.22'; $c.Style = $s
$ws.Range("E9").Value = '  -3.80%  '

# Row 10 - Toncoin
$c = $ws.Range("D10"); $s = $c.Style; $c.NumberFormat = "@"; $c.Value = '6.70'; $c.Style = $s
$ws.Range("E10").Value = '  +4.56%  '

# Row 11 - Dogecoin
$ws.Range("E11").Value = '  -4.33%  '

# Row 12 - Cardano
$c = $ws.Range("D12"); $s = $c.Style; $c.NumberFormat = "@"; $c.Value = '0.448'; $c.Style = $s
$ws.Range("E12").Value = '  -2.59%  '

# Row 13 - ShibaInu
$c = $ws.Range("D13"); $s = $c.Style; $c.NumberFormat = "@"; $c.Value = '0.0000224'; $c.Style = $s
$ws.Range("E13").Value = '  -3.66%  '

# Row 14 - Avalanche
$c = $ws.Range("D14"); $s = $c.Style; $c.NumberFormat = "@"; $c.Value = '34.01'; $c.Style = $s
$ws.Range("E14").Value = '  -1.49%  '

# Row 15 - TRON
$ws.Range("E15").Value = '  +0.76%  '

# Row 16 - WrappedliquidstakedEther2.0
$c = $ws.Range("D16"); $s = $c.Style; $c.NumberFormat = "@"; $c.Value = '3.386.50'; $c.Style = $s
$ws.Range("E16").Value = '  -3.96%  '

# Row 17 - WrappedBTC (was Polkadot, row swap w/18)
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$c = $ws.Range("D17"); $s = $c.Style; $c.NumberFormat = "@"; $c.Value = '60.773.55'; $c.Style = $s
$ws.Range("E17").Value = '  -2.67%  '

# Row 18 - Polkadot (was WrappedBTC, row swap w/17)
$ws.Range("B18").Value = 'Polkadot'
$ws.Range("C18").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$c = $ws.Range("D18"); $s = $c.Style; $c.NumberFormat = "@"; $c.Value = '6.81'; $c.Style = $s
$ws.Range("E18").Value = '  -3.04%  '

# Row 19 - WrappedEther
$c = $ws.Range("D19"); $s = $c.Style; $c.NumberFormat = "@"; $c.Value = '2.906.04'; $c.Style = $s
$ws.Range("E19").Value = '  -3.83%  '

# Row 20 - BitcoinCash
$c = $ws.Range("D20"); $s = $c.Style; $c.NumberFormat = "@"; $c.Value = '426.49'; $c.Style = $s
$ws.Range("E20").Value = '  -4.98%  '

# Row 21 - Chainlink
$c = $ws.Range("D21"); $s = $c.Style; $c.NumberFormat = "@"; $c.Value = '13.62'; $c.Style = $s
$ws.Range("E21").Value = '  -4.24%  '

# Row 22 - Polygon
$c = $ws.Range("D22"); $s = $c.Style; $c.NumberFormat = "@"; $c.Value = '0.670'; $c.Style = $s
$ws.Range("E22").Value = '  -3.09%  '

# Row 23 - Uniswap
$c = $ws.Range("D23"); $s = $c.Style; $c.NumberFormat = "@"; $c.Value = '7.08'; $c.Style = $s
$ws.Range("E23").Value = '  -4.80%  '

# Row 24 - Litecoin
$c = $ws.Range("D24"); $s = $c.Style; $c.NumberFormat = "@"; $c.Value = '80.45'; $c.Style = $s
$ws.Range("E24").Value = '  -2.29%  '

# Row 25 - RenderToken
$c = $ws.Range("D25"); $s = $c.Style; $c.NumberFormat = "@"; $c.Value = '11.00'; $c.Style = $s
$ws.Range("E25").Value = '  +2.14%  '

# Row 26 - Fetch.AI
$c = $ws.Range("D26"); $s = $c.Style; $c.NumberFormat = "@"; $c.Value = '2.20'; $c.Style = $s
$ws.Range("E26").Value = '  -1.56%  '

# Row 27 - InternetComputer(DFINITY)
$c = $ws.Range("D27"); $s = $c.Style; $c.NumberFormat = "@"; $c.Value = '11.88'; $c.Style = $s
$ws.Range("E27").Value = '  -1.03%  '

# Row 28 - Dai
$ws.Range("E28").Value = '  -0.07%  '

# Row 29 - FirstDigitalUSD
$c = $ws.Range("D29"); $s = $c.Style; $c.NumberFormat = "@"; $c.Value = '1.00'; $c.Style = $s
$ws.Range("E29").Value = '  -0.14%  '

# Row 30 - ImmutableX
$c = $ws.Range("D30"); $s = $c.Style; $c.NumberFormat = "@"; $c.Value = '2.19'; $c.Style = $s
$ws.Range("E30").Value = '  +2.22%  '

# Row 31 - NEARProtocol
$c = $ws.Range("D31"); $s = $c.Style; $c.NumberFormat = "@"; $c.Value = '7.20'; $c.Style = $s
$ws.Range("E31").Value = '  +0.38%  '

# Row 32 - PancakeSwap
$ws.Range("E32").Value = '  -3.50%  '

# Row 33 - EthereumClassic
$c = $ws.Range("D33"); $s = $c.Style; $c.NumberFormat = "@"; $c.Value = '26.47'; $c.Style = $s
$ws.Range("E33").Value = '  -4.04%  '

# Row 34 - Hedera
$c = $ws.Range("D34"); $s = $c.Style; $c.NumberFormat = "@"; $c.Value = '0.106'; $c.Style = $s
$ws.Range("E34").Value = '  -3.47%  '

# Row 35 - PEPE
$c = $ws.Range("D35"); $s = $c.Style; $c.NumberFormat = "@"; $c.Value = '0.0₃0834'; $c.Style = $s
$ws.Range("E35").Value = '  -1.66%  '

# Row 36 - Mantle
$c = $ws.Range("D36"); $s = $c.Style; $c.NumberFormat = "@"; $c.Value = '1.00'; $c.Style = $s
$ws.Range("E36").Value = '  -2.32%  '

# Row 37 - Filecoin
$c = $ws.Range("D37"); $s = $c.Style; $c.NumberFormat = "@"; $c.Value = '5.66'; $c.Style = $s
$ws.Range("E37").Value = '  -3.19%  '

# Row 38 - OKB (was Stacks, row swap w/39)
$ws.Range("B38").Value = 'OKB'
$ws.Range("C38").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$c = $ws.Range("D38"); $s = $c.Style; $c.NumberFormat = "@"; $c.Value = '49.32'; $c.Style = $s
$ws.Range("E38").Value = '  -1.60%  '

# Row 39 - Stacks (was OKB, row swap w/38)
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c = $ws.Range("D39"); $s = $c.Style; $c.NumberFormat = "@"; $c.Value = '2.03'; $c.Style = $s
$ws.Range("E39").Value = '  -2.47%  '

# Row 40 - dogwifhat
$c = $ws.Range("D40"); $s = $c.Style; $c.NumberFormat = "@"; $c.Value = '2.93'; $c.Style = $s
$ws.Range("E40").Value = '  -2.95%  '

# Row 41 - Kaspa
$ws.Range("E41").Value = '  +0.39%  '

# Row 42 - Cosmos
$c = $ws.Range("D42"); $s = $c.Style; $c.NumberFormat = "@"; $c.Value = '8.71'; $c.Style = $s
$ws.Range("E42").Value = '  -3.66%  '

# Row 43 - TheGraph
$c = $ws.Range("D43"); $s = $c.Style; $c.NumberFormat = "@"; $c.Value = '0.290'; $c.Style = $s
$ws.Range("E43").Value = '  +2.46%  '

# Row 44 - Arweave
$c = $ws.Range("D44"); $s = $c.Style; $c.NumberFormat = "@"; $c.Value = '41.87'; $c.Style = $s
$ws.Range("E44").Value = '  +4.27%  '

# Row 45 - VeChain
$c = $ws.Range("D45"); $s = $c.Style; $c.NumberFormat = "@"; $c.Value = '0.0346'; $c.Style = $s
$ws.Range("E45").Value = '  -1.92%  '

# Row 46 - Bittensor
$c = $ws.Range("D46"); $s = $c.Style; $c.NumberFormat = "@"; $c.Value = '372.86'; $c.Style = $s
$ws.Range("E46").Value = '  -4.95%  '

# Row 47 - Monero
$c = $ws.Range("D47"); $s = $c.Style; $c.NumberFormat = "@"; $c.Value = '133.76'; $c.Style = $s
$ws.Range("E47").Value = '  -0.12%  '

# Row 48 - Maker
$c = $ws.Range("D48"); $s = $c.Style; $c.NumberFormat = "@"; $c.Value = '2.660.20'; $c.Style = $s

# Row 50 - InjectiveProtocol
$c = $ws.Range("D50"); $s = $c.Style; $c.NumberFormat = "@"; $c.Value = '24.87'; $c.Style = $s
$ws.Range("E50").Value = '  +4.78%  '

# Row 51 - Stellar
$ws.Range("E51").Value = '  -1.33%  '